# Add a new worksheet "Long Rambly Name" after the existing "Ruins" sheet
# and populate it with a small table (6 rows x 3 columns), reusing the
# same cell style already used on the other sheets.

$wb = $excel.ActiveWorkbook

$settlements = $wb.Worksheets.Item("Settlements")
$ruins = $wb.Worksheets.Item("Ruins")

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ruins)
$newSheet.Name = "Long Rambly Name"

# Data for the new sheet: column A is a row number, columns B/C are text.
$data = @(
    @(1, "eheh",  "5d6+4"),
    @(2, "hello", "2d8+1"),
    @(3, "svet",  "3d2"),
    @(4, "nice",  "8d2+1"),
    @(5, "table", "100d5+1"),
    @(6, "innit", "4d10")
)

# Grab the formatting used by row 1 of "Settlements" (A1:C1) so the new
# sheet's cells end up with the same style as the rest of the workbook.
$styleSrc = $settlements.Range("A1:C1")
$styleSrc.Copy()

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $rowData = $data[$i]

    $dstRow = $newSheet.Range("A" + $row + ":C" + $row)
    $dstRow.PasteSpecial(-4122)

    $newSheet.Cells.Item($row, 1).Value = $rowData[0]
    $newSheet.Cells.Item($row, 2).Value = $rowData[1]
    $newSheet.Cells.Item($row, 3).Value = $rowData[2]
}

Write-Output "Added sheet 'Long Rambly Name' with $($data.Count) rows"
